$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2360.39455573452
$ws.Range("I2").Value = 734.394555734517

$ws.Range("B3").Value = 2352.32559555366
$ws.Range("I3").Value = 1307.32559555366

$ws.Range("B4").Value = 3170.93958214168
$ws.Range("C4").Value = 2821.15113400456
$ws.Range("I4").Value = 1529.93958214168

$ws.Range("B5").Value = 2938.06367494595
$ws.Range("E5").Value = 3307.93184266879
$ws.Range("I5").Value = 1206.06367494595

$ws.Range("B6").Value = 2614.75849661334
$ws.Range("E6").Value = 3002.2587283695
$ws.Range("F6").Value = 3181.78860470656
$ws.Range("I6").Value = 1029.75849661334

$ws.Range("B7").Value = 2238.50348973571
$ws.Range("E7").Value = 2625.72882298163
$ws.Range("I7").Value = 530.503489735707

$ws.Range("B8").Value = 2250.84188571033
$ws.Range("I8").Value = 491.841885710326

$ws.Range("B9").Value = 2251.51253617101
$ws.Range("I9").Value = 492.51253617101

$ws.Range("B10").Value = 2024.69787883702
$ws.Range("I10").Value = 213.69787883702

$ws.Range("B11").Value = 1853.53740933331
$ws.Range("I11").Value = 347.537409333307

$ws.Range("B12").Value = 2014.47765952703
$ws.Range("I12").Value = 375.477659527029

$ws.Range("B13").Value = 2082.13414466569
$ws.Range("I13").Value = 264.134144665692

$ws.Range("B14").Value = 2379.62098943604
$ws.Range("I14").Value = 541.620989436036

$ws.Range("B15").Value = 2375.27561271413
$ws.Range("I15").Value = 526.275612714135

$ws.Range("B16").Value = 3085.37024305792
$ws.Range("C16").Value = 2465.74140743452
$ws.Range("I16").Value = 94.3702430579224

$ws.Range("B17").Value = 2922.68121387256
$ws.Range("E17").Value = 3459.66053386259
$ws.Range("I17").Value = 125.681213872556

$ws.Range("B18").Value = 2654.48048423851
$ws.Range("E18").Value = 3196.43259445353
$ws.Range("I18").Value = 145.480484238507

$ws.Range("B19").Value = 2275.20377989122
$ws.Range("E19").Value = 2780.7649372921
$ws.Range("I19").Value = 9.20377989122335

$ws.Range("B20").Value = 2225.56607818956
$ws.Range("I20").Value = 85.5660781895554

$ws.Range("B21").Value = 2234.12562072625
$ws.Range("I21").Value = 125.125620726247

$ws.Range("B22").Value = 2040.98833102556
$ws.Range("I22").Value = 7.98833102555818

$ws.Range("B23").Value = 1892.94888338909
$ws.Range("I23").Value = 46.9488833890862

$ws.Range("B24").Value = 2007.09868395807
$ws.Range("I24").Value = -47.9013160419336

$ws.Range("B25").Value = 2081.45596042113
$ws.Range("I25").Value = -287.544039578875

$ws.Range("B26").Value = 2374.73095176211
$ws.Range("I26").Value = -155.269048237889

$ws.Range("B27").Value = 2325.42763219352
$ws.Range("I27").Value = -117.572367806476

$ws.Range("B28").Value = 3087.36353973135
$ws.Range("C28").Value = 2448.47204028956
$ws.Range("I28").Value = -43.6364602686485

$ws.Range("B29").Value = 2834.58864534412
$ws.Range("E29").Value = 3580.16001318352
$ws.Range("I29").Value = 309.588645344117

$ws.Range("B30").Value = 2637.18996904426
$ws.Range("E30").Value = 3240.08338226704
$ws.Range("I30").Value = 134.189969044259

$ws.Range("B31").Value = 2303.37397941122
$ws.Range("E31").Value = 2901.81894926931
$ws.Range("I31").Value = -107.626020588782

$ws.Range("B32").Value = 2243.79823076496
$ws.Range("I32").Value = 18.7982307649622

$ws.Range("B33").Value = 2237.96058241431
$ws.Range("I33").Value = -17.0394175856936

$ws.Range("B34").Value = 2057.2892337754
$ws.Range("I34").Value = 111.289233775401

$ws.Range("B35").Value = 1898.93375166217
$ws.Range("I35").Value = 215.933751662174

$ws.Range("B36").Value = 2016.55427857588
$ws.Range("I36").Value = 284.554278575883

$ws.Range("B37").Value = 2057.19702818581
$ws.Range("I37").Value = 870.197028185814

$ws.Range("B38").Value = 2338.58315343903
$ws.Range("I38").Value = 1194.58315343903

$ws.Range("B39").Value = 2333.21289648938
$ws.Range("I39").Value = 126.212896489378

$ws.Range("B40").Value = 3083.30892494582
$ws.Range("I40").Value = -33.691075054182
$ws.Range("J40").Value = "Increase"

$ws.Range("B41").Value = 2720.44158009397
$ws.Range("I41").Value = -285.558419906026

$ws.Range("B42").Value = 2657.5648211793
$ws.Range("E42").Value = 3342.17526555114
$ws.Range("I42").Value = -164.435178820697

$ws.Range("B43").Value = 2359.02527019689
$ws.Range("E43").Value = 3013.82683868839
$ws.Range("I43").Value = -169.974729803109

$ws.Range("B44").Value = 2245.03757081197
$ws.Range("I44").Value = -374.962429188034

$ws.Range("B45").Value = 2220.69897041823
$ws.Range("I45").Value = -397.301029581767

$ws.Range("B46").Value = 2061.05175497573
$ws.Range("I46").Value = -488.948245024274

$ws.Range("B47").Value = 1911.22120405076
$ws.Range("I47").Value = -631.778795949243

$ws.Range("B48").Value = 2059.46726197347
$ws.Range("I48").Value = -691.532738026527

$ws.Range("B49").Value = 2055.77357235953
$ws.Range("I49").Value = -612.226427640468

Write-Output "edit complete"
